# Updates cryptos list (prices in column D, 1h volume/change in column E)
# to the latest scrape. Two rows (45/46 and 49/50) also swap coin
# identity (name + link) because the ranking order changed.
#
# Note: several "Price" values are plain decimals (e.g. "251.80"); a
# leading apostrophe is used so Excel keeps them as text (matching the
# existing text-typed Price column) instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.940.94"
$ws.Range("E2").Value = "  +4.25%  "

$ws.Range("D3").Value = "2.281.05"
$ws.Range("E3").Value = "  +4.58%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'251.80"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  +3.78%  "

$ws.Range("D7").Value = "'72.16"
$ws.Range("E7").Value = "  +8.60%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("E9").Value = "  +15.38%  "

$ws.Range("D10").Value = "'39.24"
$ws.Range("E10").Value = "  +8.31%  "

$ws.Range("D11").Value = "'60.09"
$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("D12").Value = "'0.0969"
$ws.Range("E12").Value = "  +4.20%  "

$ws.Range("D13").Value = "'7.46"
$ws.Range("E13").Value = "  +7.90%  "

$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("D15").Value = "2.623.08"
$ws.Range("E15").Value = "  +4.61%  "

$ws.Range("D16").Value = "'14.99"
$ws.Range("E16").Value = "  +4.17%  "

$ws.Range("D17").Value = "'0.889"
$ws.Range("E17").Value = "  +4.00%  "

$ws.Range("D18").Value = "2.279.83"
$ws.Range("E18").Value = "  +2.71%  "

$ws.Range("D19").Value = "42.880.20"
$ws.Range("E19").Value = "  +4.16%  "

$ws.Range("E20").Value = "  +7.83%  "

$ws.Range("E21").Value = "  +3.94%  "

$ws.Range("D22").Value = "'73.53"
$ws.Range("E22").Value = "  +2.58%  "

$ws.Range("D23").Value = "'234.46"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("D24").Value = "'2.19"
$ws.Range("E24").Value = "  +7.54%  "

$ws.Range("D25").Value = "'4.04"
$ws.Range("E25").Value = "  +7.35%  "

$ws.Range("D26").Value = "'11.53"
$ws.Range("E26").Value = "  +1.48%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "'2.46"
$ws.Range("E28").Value = "  +1.43%  "

$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").Value = "'167.88"
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").Value = "'21.10"
$ws.Range("E32").Value = "  +4.25%  "

$ws.Range("D33").Value = "'6.59"
$ws.Range("E33").Value = "  +14.28%  "

$ws.Range("D34").Value = "'0.128"
$ws.Range("E34").Value = "  +4.13%  "

$ws.Range("D35").Value = "'32.09"
$ws.Range("E35").Value = "  +31.37%  "

$ws.Range("D36").Value = "'0.0804"
$ws.Range("E36").Value = "  +9.33%  "

$ws.Range("D37").Value = "'0.127"
$ws.Range("E37").Value = "  +4.17%  "

$ws.Range("D38").Value = "'4.50"
$ws.Range("E38").Value = "  +13.82%  "

$ws.Range("E39").Value = "  +5.74%  "

$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").Value = "'13.45"
$ws.Range("E41").Value = "  +18.32%  "

$ws.Range("E42").Value = "  +5.36%  "

$ws.Range("D43").Value = "'5.85"
$ws.Range("E43").Value = "  +6.57%  "

$ws.Range("E44").Value = "  +11.17%  "

$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").Value = "'62.69"
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.20"
$ws.Range("E46").Value = "  +7.67%  "

$ws.Range("D47").Value = "'5.01"
$ws.Range("E47").Value = "  -5.28%  "

$ws.Range("E48").Value = "  +2.47%  "

$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.19"
$ws.Range("E50").Value = "  +3.73%  "

$ws.Range("E51").Value = "  +4.00%  "
